# Fruta / hortaliza, semanal
# Sort the data rows (2-5) of the sheet by the "Fecha" column (D) ascending.
# This reorders the Volumen/Precio mínimo/máximo/promedio/Precio $/Kg columns
# (J, K, L, M, P) along with the date (D) while the remaining columns stay
# the same because their values are identical across all rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44277
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 11000
$ws.Range("M2").Value = 11000
$ws.Range("P2").Value = 550

$ws.Range("D3").Value = 44280
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 500

$ws.Range("D5").Value = 44291
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 11000
$ws.Range("M5").Value = 11000
$ws.Range("P5").Value = 550
